# Reducing memory consumption:
# Drop the intermediate / no-longer-needed computed columns
# (the old index column A, and peak_runoff_rate / infiltration / evaporation)
# and keep only "runoff" and "Destore-Imperv".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove peak_runoff_rate / infiltration / evaporation columns (old C:E).
# This shifts the old Destore-Imperv column (F) left into column C.
$ws.Range("C1:E1").EntireColumn.Delete() | Out-Null

# Remove the old numeric index column (A).
# This shifts runoff (old B) into A and Destore-Imperv (old C, after the
# previous delete) into B.
$ws.Range("A1").EntireColumn.Delete() | Out-Null

# The header row formatting (bold, centered, bordered) that used to live on
# the deleted header cells still lingers on the now-empty C1:E1 cells;
# reproduce that by copying the format from the surviving header cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection where it ended up after working in the sheet.
$ws.Range("F6").Select() | Out-Null
